$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.823.74"
$ws.Cells.Item(2, 5).Value = "  +0.02%  "

$ws.Cells.Item(3, 4).Value = "1.631.44"
$ws.Cells.Item(3, 5).Value = "  -0.62%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.42%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "214.28"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.71%  "

$ws.Cells.Item(6, 5).Value = "  -0.37%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.46%  "

$ws.Cells.Item(8, 5).Value = "  -1.27%  "

$ws.Cells.Item(9, 5).Value = "  -0.92%  "

$ws.Cells.Item(10, 5).Value = "  -0.04%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.0790"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.83%  "

$ws.Cells.Item(12, 5).Value = "  +0.02%  "

$ws.Cells.Item(13, 4).Value = "1.857.07"
$ws.Cells.Item(13, 5).Value = "  -0.16%  "

$ws.Cells.Item(14, 4).Value = "1.636.42"
$ws.Cells.Item(14, 5).Value = "  -0.70%  "

$ws.Cells.Item(15, 5).Value = "  -0.67%  "

$ws.Cells.Item(16, 5).Value = "  -0.92%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "62.79"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.27%  "

$ws.Cells.Item(18, 4).Value = "25.814.40"
$ws.Cells.Item(18, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.47%  "

$ws.Cells.Item(20, 5).Value = "  -0.07%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "191.32"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.89%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "9.92"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.17%  "

$ws.Cells.Item(23, 5).Value = "  +0.20%  "

$ws.Cells.Item(24, 5).Value = "  +2.14%  "

$ws.Cells.Item(25, 5).Value = "  -0.51%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "142.37"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.73%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.98%  "

$ws.Cells.Item(28, 5).Value = "  -0.60%  "

$ws.Cells.Item(29, 5).Value = "  -0.23%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.Style = "Normal"

$ws.Cells.Item(31, 5).Value = "  +0.63%  "

$ws.Cells.Item(32, 5).Value = "  -0.60%  "

$ws.Cells.Item(33, 5).Value = "  -1.15%  "

$ws.Cells.Item(34, 5).Value = "  -0.29%  "

$ws.Cells.Item(35, 5).Value = "  +0.00%  "

$ws.Cells.Item(36, 5).Value = "  +0.47%  "

$ws.Cells.Item(37, 4).Value = "1.141.70"
$ws.Cells.Item(37, 5).Value = "  +2.62%  "

$ws.Cells.Item(38, 5).Value = "  -0.53%  "

$ws.Cells.Item(39, 5).Value = "  -2.35%  "

$ws.Cells.Item(40, 5).Value = "  -0.45%  "

$ws.Cells.Item(41, 5).Value = "  -0.61%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.51"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.28%  "

$ws.Cells.Item(43, 5).Value = "  +0.36%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "100.62"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.49%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.802"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.03%  "

$ws.Cells.Item(46, 4).Value = "1.766.79"
$ws.Cells.Item(46, 5).Value = "  +0.24%  "

$ws.Cells.Item(47, 5).Value = "  +0.10%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "55.37"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +0.29%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.47"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +6.63%  "

$ws.Cells.Item(50, 5).Value = "  +1.94%  "

$ws.Cells.Item(51, 5).Value = "  -0.71%  "
